$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update parentCriterion values in column J: "outcomes" -> "criteria" ---
# (rows 3, 6, 9 hold the second-level groups whose parent is the top "criteria" row)
$ws.Range("J3").Value = "criteria"
$ws.Range("J6").Value = "criteria"
$ws.Range("J9").Value = "criteria"

# --- Add new column K: "isLeaf" ---
# Copy the header formatting from J1 (grey fill + wrap) onto K1 so the new
# header cell matches the existing header style instead of minting a new one.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K1").Value = "isLeaf"

$ws.Range("K2").Value = $false
$ws.Range("K3").Value = $false
$ws.Range("K4").Value = $true
$ws.Range("K5").Value = $true
$ws.Range("K6").Value = $false
$ws.Range("K7").Value = $true
$ws.Range("K8").Value = $true
$ws.Range("K9").Value = $false
$ws.Range("K10").Value = $true
$ws.Range("K11").Value = $true

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 26.0
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 31.666666666666668
$ws.Columns.Item(10).ColumnWidth = 15.833333333333334

# --- Row heights (Excel re-wraps text after the column-width changes above) ---
# Rows 3, 6 and 9 are short (non-wrapping) rows that go back to the sheet's
# default height once the content no longer wraps - AutoFit reverts them to
# the plain, non-custom default height.
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(9).AutoFit()

# Rows 4, 5, 7, 8, 10 and 11 wrap onto fewer lines now that the text moved
# into the wider B/C columns, so their height shrinks accordingly.
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 57.6
$ws.Rows.Item(7).RowHeight = 28.8
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(10).RowHeight = 43.2
$ws.Rows.Item(11).RowHeight = 28.8

# --- Selection change ---
$ws.Range("K10").Select()
